# Auto-generated edit script: apply 2024-04-19 data updates
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 2137
$ws.Range("K3").Value = 2070
$ws.Range("J4").Value = 1807
$ws.Range("K4").Value = 432
$ws.Range("I5").Value = 724
$ws.Range("K5").Value = 138
$ws.Range("K6").Value = 2622
$ws.Range("I7").Value = 26239
$ws.Range("J7").Value = 29278
$ws.Range("K7").Value = 7399

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 144
$ws.Range("K3").Value = 142
$ws.Range("K7").Value = 498

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 64
$ws.Range("K3").Value = 50
$ws.Range("K7").Value = 158

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 83
$ws.Range("K6").Value = 77
$ws.Range("K7").Value = 292

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 41
$ws.Range("K7").Value = 115

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K3").Value = 81
$ws.Range("K7").Value = 239

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 52
$ws.Range("K3").Value = 44
$ws.Range("K6").Value = 75
$ws.Range("K7").Value = 181

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K6").Value = 41
$ws.Range("K7").Value = 136

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 55
$ws.Range("K5").Value = 11
$ws.Range("K8").Value = 498
$ws.Range("K10").Value = 42
$ws.Range("K11").Value = 160
$ws.Range("K18").Value = 51
$ws.Range("K19").Value = 208
$ws.Range("K20").Value = 157
$ws.Range("K21").Value = 21
$ws.Range("K22").Value = 23
$ws.Range("K29").Value = 372
$ws.Range("K31").Value = 83
$ws.Range("K32").Value = 12
$ws.Range("K33").Value = 292
$ws.Range("K36").Value = 87
$ws.Range("K37").Value = 239
$ws.Range("K42").Value = 256
$ws.Range("K49").Value = 51
$ws.Range("K50").Value = 46
$ws.Range("K52").Value = 198
$ws.Range("I63").Value = 198
$ws.Range("J63").Value = 96
$ws.Range("K63").Value = 24
$ws.Range("K65").Value = 181
$ws.Range("K67").Value = 282
$ws.Range("K68").Value = 19
$ws.Range("K72").Value = 35
$ws.Range("K75").Value = 31
$ws.Range("K79").Value = 198
$ws.Range("K81").Value = 7
$ws.Range("K83").Value = 158
$ws.Range("K85").Value = 364
$ws.Range("K86").Value = 50
$ws.Range("K89").Value = 101
$ws.Range("K90").Value = 63
$ws.Range("J93").Value = 121
$ws.Range("K93").Value = 34
$ws.Range("K95").Value = 115
$ws.Range("K96").Value = 104
$ws.Range("K99").Value = 136
$ws.Range("K100").Value = 11
$ws.Range("I101").Value = 26239
$ws.Range("J101").Value = 29278
$ws.Range("K101").Value = 7399

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K6").Value = 34
$ws.Range("K7").Value = 83

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 82
$ws.Range("K7").Value = 282

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K2").Value = 5
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 97
$ws.Range("K3").Value = 125
$ws.Range("K5").Value = 10
$ws.Range("K6").Value = 120
$ws.Range("K7").Value = 372

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K3").Value = 58
$ws.Range("K7").Value = 208

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 61
$ws.Range("K7").Value = 256

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 42

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K6").Value = 50
$ws.Range("K7").Value = 104

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 21

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K3").Value = 72
$ws.Range("K4").Value = 14
$ws.Range("K7").Value = 198

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K3").Value = 45
$ws.Range("K7").Value = 157

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K2").Value = 16
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K2").Value = 34
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("K2").Value = 11
$ws.Range("J4").Value = 10
$ws.Range("J7").Value = 121
$ws.Range("K7").Value = 34

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("K2").Value = 3
$ws.Range("K7").Value = 11

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 46

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K4").Value = 7
$ws.Range("K6").Value = 66
$ws.Range("K7").Value = 160

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K2").Value = 16
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 12

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K2").Value = 21
$ws.Range("K7").Value = 101

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 11

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K2").Value = 11
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("K3").Value = 7
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 63

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("K2").Value = 6
$ws.Range("K7").Value = 19

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K3").Value = 122
$ws.Range("K6").Value = 88
$ws.Range("K7").Value = 364

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("K2").Value = 11
$ws.Range("K7").Value = 23

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K3").Value = 44
$ws.Range("K6").Value = 86
$ws.Range("K7").Value = 198

$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Range("K3").Value = 1
$ws.Range("K6").Value = 4
$ws.Range("K7").Value = 7
